$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose value changed (or are newly populated). Stored as a hashtable
# mapping 1-based column index -> new text. We force text/shared-string storage
# (NumberFormat "@" then ClearFormats afterwards) so that numeric-looking values
# such as "15" or "1234" stay text instead of turning into Excel numbers,
# matching how the source data is stored.
$updates = @{
    6 = "Universitari grau superior"   # F2
    7 = "15"   # G2
    10 = "Militar"   # J2
    322 = "18"   # LJ2
    323 = "24"   # LK2
    324 = "25"   # LL2
    325 = "23.08"   # LM2
    326 = "48.08"   # LN2
    327 = "38"   # LO2
    328 = "15"   # LP2
    329 = "8"   # LQ2
    330 = "40"   # LR2
    331 = "43"   # LS2
    332 = "40"   # LT2
    333 = "48"   # LU2
    334 = "96"   # LV2
    347 = "6"   # MI2
    348 = "11-18"   # MJ2
    349 = "7"   # MK2
    350 = "4"   # ML2
    351 = "8"   # MM2
    352 = "60-71"   # MN2
    353 = "11"   # MO2
    354 = "19"   # MP2
    355 = "9"   # MQ2
    356 = "44763"   # MR2
    357 = "0.42053"   # MS2
    358 = "2"   # MT2
    359 = "0"   # MU2
    360 = "12"   # MV2
    361 = "46952"   # MW2
    362 = "0.70167"   # MX2
    363 = "5"   # MY2
    364 = "0"   # MZ2
    365 = "12"   # NA2
    366 = "37549"   # NB2
    367 = "104720"   # NC2
    368 = "9"   # ND2
    369 = "0"   # NE2
    370 = "12"   # NF2
    371 = "35"   # NG2
    372 = "2"   # NH2
    373 = "246102"   # NI2
    374 = "0.10739"   # NJ2
    375 = "133615"   # NK2
    376 = "30"   # NL2
    377 = "7"   # NM2
    378 = "268804"   # NN2
    379 = "0.08108"   # NO2
    380 = "112080"   # NP2
    381 = "55"   # NQ2
    382 = "33"   # NR2
    383 = "308798"   # NS2
    384 = "0.13869"   # NT2
    385 = "0.91174"   # NU2
    386 = "31"   # NV2
    387 = "10"   # NW2
    388 = "296097"   # NX2
    389 = "0.15355"   # NY2
    390 = "105427"   # NZ2
    391 = "31"   # OA2
    392 = "18"   # OB2
    393 = "368781"   # OC2
    394 = "0.14723"   # OD2
    395 = "0.91966"   # OE2
    396 = "5"   # OF2
    397 = "2"   # OG2
    398 = "0.70167"   # OH2
    399 = "76201"   # OI2
    400 = "16"   # OJ2
    401 = "10"   # OK2
    402 = "15"   # OL2
    403 = "11"   # OM2
    404 = "4"   # ON2
    405 = "19-28"   # OO2
    406 = "8"   # OP2
    433 = "00:00:00"   # PQ2
    434 = "1"   # PR2
    435 = "0"   # PS2
    436 = "1"   # PT2
    441 = "00:00:00"   # PY2
    443 = "0"   # QA2
    444 = "1"   # QB2
    445 = "1"   # QC2
    446 = "1"   # QD2
    447 = "444"   # QE2
    448 = "444"   # QF2
    449 = "444"   # QG2
    450 = "444"   # QH2
    452 = "1"   # QJ2
    453 = "11"   # QK2
    454 = "25"   # QL2
    455 = "13.46"   # QM2
    456 = "38.46"   # QN2
    457 = "8"   # QO2
    458 = "8"   # QP2
    459 = "10"   # QQ2
    460 = "30"   # QR2
    461 = "20"   # QS2
    462 = "15"   # QT2
    463 = "42"   # QU2
    464 = "54"   # QV2
    474 = "4"   # RF2
    475 = "6"   # RG2
    476 = "11-18"   # RH2
    477 = "7"   # RI2
    478 = "0"   # RJ2
    479 = "0"   # RK2
    480 = "0"   # RL2
    481 = "0"   # RM2
    482 = "3"   # RN2
    483 = "0"   # RO2
    528 = "0"   # TH2
    529 = "0"   # TI2
    530 = "0"   # TJ2
    531 = "0"   # TK2
    532 = "0"   # TL2
    533 = "0"   # TM2
    534 = "0"   # TN2
    557 = "00:00:00"   # UK2
    559 = "0"   # UM2
    569 = "00:00:00"   # UW2
    571 = "0"   # UY2
    580 = "3"   # VH2
    581 = "3"   # VI2
    582 = "8.33"   # VJ2
    583 = "0"   # VK2
    584 = "8.33"   # VL2
    585 = "0"   # VM2
    586 = "3"   # VN2
    587 = "0"   # VO2
    588 = "10"   # VP2
    589 = "6"   # VQ2
    590 = "13"   # VR2
    591 = "7"   # VS2
    592 = "16"   # VT2
}

foreach ($col in $updates.Keys) {
    $text = $updates[$col]
    $cell = $ws.Cells.Item(2, $col)
    $cell.NumberFormat = "@"
    $cell.Value = $text
    $cell.ClearFormats()
}

# Cells that previously held a value but should become blank.
$clearCols = @(
    342,  # MD2
    343,  # ME2
    344,  # MF2
    345,  # MG2
    425,  # PI2
    426,  # PJ2
    427,  # PK2
    428,  # PL2
    437,  # PU2
    438,  # PV2
    439,  # PW2
    440,  # PX2
    466,  # QX2
    467,  # QY2
    468,  # QZ2
    469,  # RA2
    470,  # RB2
    471,  # RC2
    472,  # RD2
    473,  # RE2
    520,  # SZ2
    521,  # TA2
    522,  # TB2
    523,  # TC2
    524,  # TD2
    525,  # TE2
    526,  # TF2
    549,  # UC2
    551,  # UE2
    553,  # UG2
    555,  # UI2
)
foreach ($col in $clearCols) {
    $ws.Cells.Item(2, $col).ClearContents()
}

Write-Host "Done applying updates"
